$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows down (format
# taken from above, i.e. from the header row, is then cleaned up below).
$ws.Rows.Item(2).Insert(-4121, 0)

# The freshly inserted row inherits the header's bold/border style;
# strip that so the new row matches the plain style of the data rows.
$ws.Range("A2:O2").ClearFormats()

# Re-apply the date format used by the other rows' arrival/departure cells.
$ws.Range("D2:E2").NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(2, 1).Value = "Cristina Testoni"
$ws.Cells.Item(2, 2).Value = "Booking"

# Force the phone number to stay text (it has a leading "+") instead of
# being auto-coerced into a number, then drop the format flag we used to
# do that so the cell ends up with the default style.
$ws.Range("C2").NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "+393316017468"
$ws.Range("C2").ClearFormats()

$ws.Cells.Item(2, 4).Value = 45597
$ws.Cells.Item(2, 5).Value = 45880
$ws.Cells.Item(2, 6).Value = 283
$ws.Cells.Item(2, 7).Value = 244.34
$ws.Cells.Item(2, 8).Value = 200.75
$ws.Cells.Item(2, 9).Value = 43.59
$ws.Cells.Item(2, 10).Value = 17.84
$ws.Cells.Item(2, 11).Value = 2024
$ws.Cells.Item(2, 12).Value = 11
